$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Merge the two runs that make up the "communication methods" paragraph
#    (they are currently split only by the _GoBack bookmark) into one run,
#    and drop the bookmark from its old position.
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$merged = [string]"没必要一直占用信道资源而不让出给其它用户使用，因此这两种方式对信道的利用率都不高。"
$content = $d.Content
$content.Find.Execute($merged, $true, $false, $false, $false, $false, $true, 1, $false, $merged, 2) | Out-Null

# Locate that paragraph again (robust to any renumbering above).
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.StartsWith("使用这两种方式进行通信")) {
        $targetIndex = $i
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Insert two new (empty) paragraphs right after the first blank paragraph
#    that follows the content paragraph, then fill them with the new
#    "Statistical Time Division Multiplexing" text.
# ---------------------------------------------------------------------------
$firstBlank = $targetIndex + 1
$p = $d.Paragraphs($firstBlank).Range
$p.InsertParagraphAfter()
$p2 = $d.Paragraphs($firstBlank).Range
$p2.InsertParagraphAfter()

$headingIndex = $firstBlank + 1
$descIndex = $firstBlank + 2

$heading = $d.Paragraphs($headingIndex).Range
$heading.Collapse(1)
$heading.InsertAfter("2. 统计时分复用")

$desc = $d.Paragraphs($descIndex).Range
$desc.Collapse(1)
$desc.InsertAfter("是对时分复用的一种改进，不固定每个用户在时分复用帧中的位置，只要有数据就集中起来组成统计时分复用帧然后发送。")

# ---------------------------------------------------------------------------
# 3. Re-create the _GoBack bookmark on the (now empty) paragraph that
#    follows the two blank paragraphs coming after the new content.
# ---------------------------------------------------------------------------
$bookmarkParaIndex = $descIndex + 2
$bmRange = $d.Paragraphs($bookmarkParaIndex).Range
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output ("Done. ParaCount=" + $d.Paragraphs.Count)
